$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formulas in C3 and D3 to be negated
$ws.Range("C3").Formula = "=-6400*1.15"
$ws.Range("D3").Formula = "=-6400*1.1"

# Update the active selection shown in the worksheet view (D4 -> C5)
$ws.Range("C5").Select()
